$wb = $excel.ActiveWorkbook

# --- Update the M-column formulas on "Capella" and "FlavourArt" sheets ---
# Original formula:  =CONCATENATE($G$1,$H$1,A{r},$I$2,$I$1,C{r},$K$1)
# New formula:        =CONCATENATE($G$1,$H$1,A{r},$I$2,$I$1,C{r},$I$2,$J$1,$K$1)
# (adds the ", isSelected: false" segment before the closing " },")

$sheetsToUpdate = @(
    @{ Name = "Capella"; LastRow = 94 },
    @{ Name = "FlavourArt"; LastRow = 39 }
)

foreach ($info in $sheetsToUpdate) {
    $ws = $wb.Worksheets.Item($info.Name)
    for ($r = 1; $r -le $info.LastRow; $r++) {
        $ws.Cells.Item($r, 13).Formula = "=CONCATENATE(`$G`$1,`$H`$1,A$r,`$I`$2,`$I`$1,C$r,`$I`$2,`$J`$1,`$K`$1)"
    }
}

# --- Update active tab / tab selection / scroll position ---
# Activating "Capella" first clears its stale topLeftCell="A49" scroll position.
# Activating "FlavourArt" last makes it the active tab (tabSelected + workbook activeTab),
# and clears tabSelected from "TPA" (which had it before).
$wb.Worksheets.Item("Capella").Activate()
$wb.Worksheets.Item("FlavourArt").Activate()
